$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header name / card number changes
$ws.Range("C2").Value = "Hartmut"
# Card number is a 16-digit string; a leading apostrophe forces Excel to
# store it as text (quote-prefixed) instead of silently coercing it to a
# Number and losing precision beyond 15 significant digits.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 06.07.2024"

# Row 6: transaction
$ws.Range("B6").Value = "07.07."
$ws.Range("C6").Value = "08.07."
$ws.Range("D6").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E6").Value = "52,70-"

# Row 7: transaction
$ws.Range("B7").Value = "11.07."
$ws.Range("C7").Value = "12.07."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 58028654"
$ws.Range("E7").Value = "87,86-"

# Row 8: transaction
$ws.Range("B8").Value = "12.07."
$ws.Range("C8").Value = "13.07."
$ws.Range("D8").Value = "BURGER KING Badoberan"
$ws.Range("E8").Value = "26,12-"

# Row 9: cleared out (was a 4th transaction, now blank) - style changes to match the now-empty row pattern
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 15.07.2024"
$ws.Range("E12").Value = "166,68-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 23.07.2024"
